# Apply odds updates to Sheet1 per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Melbourne City vs Perth Glory)
$ws.Range("K2").Value  = 4.3
$ws.Range("L2").Value  = 1.38
$ws.Range("X2").Value  = 16
$ws.Range("Z2").Value  = 48
$ws.Range("AA2").Value = 190
$ws.Range("AC2").Value = 9.4
$ws.Range("AE2").Value = 85
$ws.Range("AF2").Value = 10
$ws.Range("AG2").Value = 9.800000000000001
$ws.Range("AH2").Value = 21
$ws.Range("AI2").Value = 85
$ws.Range("AJ2").Value = 16.5
$ws.Range("AK2").Value = 17.5
$ws.Range("AL2").Value = 36
$ws.Range("AM2").Value = 130
$ws.Range("AO2").Value = 110

# Row 3 (Serik Belediyespor vs Boluspor)
$ws.Range("F3").Value = 6.2

# Row 5 (AC Milan vs Verona)
$ws.Range("S5").Value = 3.25
$ws.Range("U5").Value = 1.77

# Row 6 (Erzurum BB vs Corum Belediyespor)
$ws.Range("G6").Value = 2.88
$ws.Range("H6").Value = 2.88
$ws.Range("J6").Value = 2.98
$ws.Range("K6").Value = 5.6

# Row 7 (Sunderland vs Leeds)
$ws.Range("H7").Value = 2.98

# Row 9 (Arouca vs Gil Vicente)
$ws.Range("P9").Value = 1.68
$ws.Range("Q9").Value = 2.28

# Row 10 (Casa Pia vs Guimaraes)
$ws.Range("H10").Value = 2.36
$ws.Range("P10").Value = 1.55
$ws.Range("Q10").Value = 2.48

# Row 11 (Amed Sportif Faaliyetler vs 76 Igdir Belediyespor)
$ws.Range("Q11").Value = 1.65

# Row 12 (Crystal Palace vs Tottenham)
$ws.Range("N12").Value  = 3.75
$ws.Range("AF12").Value = 15
$ws.Range("AI12").Value = 55
$ws.Range("AJ12").Value = 34

# Row 13 (Bologna vs Sassuolo)
$ws.Range("N13").Value  = 3.7
$ws.Range("X13").Value  = 14
$ws.Range("AI13").Value = 85

# Row 14 (Braga vs Benfica)
$ws.Range("P14").Value = 1.77
$ws.Range("Q14").Value = 1.92

# Row 15 (Atalanta vs Inter)
$ws.Range("G15").Value  = 3.9
$ws.Range("N15").Value  = 4.4
$ws.Range("P15").Value  = 2.16
$ws.Range("T15").Value  = 1.71
$ws.Range("U15").Value  = 2.32
$ws.Range("AC15").Value = 8.4
$ws.Range("AH15").Value = 16.5
$ws.Range("AI15").Value = 34
$ws.Range("AK15").Value = 44
$ws.Range("AL15").Value = 50

# Row 16 (Sporting Lisbon vs Rio Ave)
$ws.Range("K16").Value = 10.5
